$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Training Dashboard")
$ws2 = $wb.Worksheets.Item("Exam Dashboard")

# --- Training Dashboard: refresh "period to expire" + "last update" numbers ---
# Row 3
$ws1.Range("H3").Value = 84
$i3 = $ws1.Range("I3")
$i3.Formula = "=""16-Sep-2025"""
$i3.Copy()
$i3.PasteSpecial(-4163)   # xlPasteValues - keeps literal text, avoids date auto-parsing

# Row 4
$ws1.Range("H4").Value = -51
$i4 = $ws1.Range("I4")
$i4.Formula = "=""16-Sep-2025"""
$i4.Copy()
$i4.PasteSpecial(-4163)

# Row 5
$ws1.Range("H5").Value = 226
$i5 = $ws1.Range("I5")
$i5.Formula = "=""16-Sep-2025"""
$i5.Copy()
$i5.PasteSpecial(-4163)

$excel.CutCopyMode = $false

# --- Exam Dashboard: narrower comments column + updated remark ---
$ws2.Columns.Item(5).ColumnWidth = 14.17   # renders as stored width 15 (engine applies a 0.83 offset)
$ws2.Range("E3").Value = "date is valid"

# --- Styling: titles & header bands drop the big 14pt font in favour of a
#     shared bold white font (headers keep their navy fill / title keeps centring) ---
$ws1.Range("A1").Font.Size = 11
$ws1.Range("A1").Font.Color = 16777215
$ws2.Range("A1").Font.Size = 11
$ws2.Range("A1").Font.Color = 16777215

$ws1.Range("A2:K2").Font.Color = 16777215
$ws2.Range("A2:G2").Font.Color = 16777215
